$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.161.84"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.568.47"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'584.56"
$ws.Range("E5").Value = "  +3.00%  "
$ws.Range("D6").Value = "'147.67"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  +2.58%  "
$ws.Range("E9").Value = "  +3.40%  "
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D12").Value = "'0.356"
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").Value = "'27.41"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "3.028.51"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "63.116.69"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("D17").Value = "2.559.43"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "'11.35"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "'343.98"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("D21").Value = "'6.89"
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -3.74%  "
$ws.Range("D24").Value = "'66.83"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").Value = "2.697.28"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "'1.64"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'8.14"
$ws.Range("E28").Value = "  +11.93%  "
$ws.Range("D29").Value = "'8.47"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("E32").Value = "  +7.83%  "
$ws.Range("D33").Value = "0.0₃0827"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").Value = "'465.32"
$ws.Range("E34").Value = "  +14.33%  "
$ws.Range("E35").Value = "  +3.67%  "
$ws.Range("D36").Value = "'175.70"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("E37").Value = "  +2.40%  "
$ws.Range("D38").Value = "'19.21"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").Value = "'4.57"
$ws.Range("E39").Value = "  +4.89%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'1.75"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'151.60"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D44").Value = "'3.83"
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("D45").Value = "'21.01"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("E46").Value = "  +5.86%  "
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("D48").Value = "'0.0976"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("E51").Value = "  -0.03%  "
